$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.7164481774827038
$ws.Range("C2").Value = 0.1654998379516996
$ws.Range("D2").Value = 0.1997473665255285
$ws.Range("E2").Value = 0.1554444581064871
$ws.Range("F2").Value = 1.203585321672342
$ws.Range("H2").Value = 0.07973214163530429
$ws.Range("I2").Value = 0.5308228575646829
$ws.Range("J2").Value = 0.1612949153979244
$ws.Range("M2").Value = 0.330032997183082
$ws.Range("N2").Value = 1.185044160131291
$ws.Range("O2").Value = 2.763206647604108

$ws.Range("B3").Value = 0.6398855889804906
$ws.Range("C3").Value = 0.1453888957765912
$ws.Range("D3").Value = 0.1980223503766538
$ws.Range("E3").Value = 0.15520883428999
$ws.Range("F3").Value = 1.202166773392321
$ws.Range("H3").Value = 0.07973214163530429
$ws.Range("I3").Value = 0.5359569668228197
$ws.Range("J3").Value = 0.1618627829399131
$ws.Range("M3").Value = 0.3090634346379275
$ws.Range("N3").Value = 1.19461734461435
$ws.Range("O3").Value = 2.763768563917949

$ws.Range("B4").Value = 0.5928852870376318
$ws.Range("C4").Value = 0.1330184562545185
$ws.Range("D4").Value = 0.1970311934780113
$ws.Range("E4").Value = 0.1551299897062073
$ws.Range("F4").Value = 1.202008855747827
$ws.Range("H4").Value = 0.07973214163530429
$ws.Range("I4").Value = 0.5394442511005018
$ws.Range("J4").Value = 0.162290088693787
$ws.Range("M4").Value = 0.2962820982477794
$ws.Range("N4").Value = 1.200953084107304
$ws.Range("O4").Value = 2.765888525533711

$ws.Range("B5").Value = 0.5737358378647457
$ws.Range("C5").Value = 0.1279720694116406
$ws.Range("D5").Value = 0.1966444599926334
$ws.Range("E5").Value = 0.1551144391353354
$ws.Range("F5").Value = 1.2021238025104
$ws.Range("H5").Value = 0.07973214163530429
$ws.Range("I5").Value = 0.5409494834176556
$ws.Range("J5").Value = 0.1624839963423703
$ws.Range("M5").Value = 0.2910975748828335
$ws.Range("N5").Value = 1.203650151577357
$ws.Range("O5").Value = 2.767198470183587

$ws.Range("B6").Value = 0.570556336013226
$ws.Range("C6").Value = 0.1271338063318979
$ws.Range("D6").Value = 0.1965812823202739
$ws.Range("E6").Value = 0.1551128590833315
$ws.Range("F6").Value = 1.202153718327594
$ws.Range("H6").Value = 0.07973214163530429
$ws.Range("I6").Value = 0.5412045052603069
$ws.Range("J6").Value = 0.1625173891744183
$ws.Range("M6").Value = 0.2902381462740706
$ws.Range("N6").Value = 1.204104958324194
$ws.Range("O6").Value = 2.767442917219881

$ws.Range("B7").Value = 0.5926270149048207
$ws.Range("C7").Value = 0.1329504201606824
$ws.Range("D7").Value = 0.1970259082371513
$ws.Range("E7").Value = 0.1551297128223759
$ws.Range("F7").Value = 1.202009679997097
$ws.Range("H7").Value = 0.07973214163530429
$ws.Range("I7").Value = 0.5394642106384318
$ws.Range("J7").Value = 0.1622926237176721
$ws.Range("M7").Value = 0.2962120804503812
$ws.Range("N7").Value = 1.20098899112795
$ws.Range("O7").Value = 2.765904386270563

$ws.Range("B8").Value = 0.6900481235519464
$ws.Range("C8").Value = 0.1585703508069685
$ws.Range("D8").Value = 0.1991385061599189
$ws.Range("E8").Value = 0.1553495672340262
$ws.Range("F8").Value = 1.202948197420781
$ws.Range("H8").Value = 0.07973214163530429
$ws.Range("I8").Value = 0.53252353583445
$ws.Range("J8").Value = 0.161474398175578
$ws.Range("M8").Value = 0.3227833566020095
$ws.Range("N8").Value = 1.188250040053227
$ws.Range("O8").Value = 2.763031829793505

$ws.Range("B9").Value = 0.8811228956523678
$ws.Range("C9").Value = 0.2086254696965852
$ws.Range("D9").Value = 0.2038182534638509
$ws.Range("E9").Value = 0.1563020753515723
$ws.Range("F9").Value = 1.210448915933384
$ws.Range("H9").Value = 0.07973214163530429
$ws.Range("I9").Value = 0.5215741899957145
$ws.Range("J9").Value = 0.1604937151416479
$ws.Range("M9").Value = 0.3756244571425356
$ws.Range("N9").Value = 1.16689719385699
$ws.Range("O9").Value = 2.771498834094132

$ws.Range("B10").Value = 1.021481722438295
$ws.Range("C10").Value = 0.2452792984621226
$ws.Range("D10").Value = 0.207580741090382
$ws.Range("E10").Value = 0.1573187209383491
$ws.Range("F10").Value = 1.219416372708011
$ws.Range("H10").Value = 0.07973214163530429
$ws.Range("I10").Value = 0.5151576068202992
$ws.Range("J10").Value = 0.1601536291777492
$ws.Range("M10").Value = 0.4148835828581809
$ws.Range("N10").Value = 1.153415967744948
$ws.Range("O10").Value = 2.78634428363236

$ws.Range("B11").Value = 1.085321043901388
$ws.Range("C11").Value = 0.2619260499059521
$ws.Range("D11").Value = 0.2093621935758563
$ws.Range("E11").Value = 0.157849823266865
$ws.Range("F11").Value = 1.224247958757971
$ws.Range("H11").Value = 0.07973214163530429
$ws.Range("I11").Value = 0.5125932227040941
$ws.Range("J11").Value = 0.1600815571731289
$ws.Range("M11").Value = 0.4328362805450681
$ws.Range("N11").Value = 1.147761183617654
$ws.Range("O11").Value = 2.794977244958147

$ws.Range("B12").Value = 1.109492791744287
$ws.Range("C12").Value = 0.2682256051931518
$ws.Range("D12").Value = 0.2100467691945198
$ws.Range("E12").Value = 0.1580607846934825
$ws.Range("F12").Value = 1.226185793187639
$ws.Range("H12").Value = 0.07973214163530429
$ws.Range("I12").Value = 0.5116732397632902
$ws.Range("J12").Value = 0.160066148413307
$ws.Range("M12").Value = 0.4396476769410071
$ws.Range("N12").Value = 1.145688519015835
$ws.Range("O12").Value = 2.798517045200356

$ws.Range("B13").Value = 1.104287119218384
$ws.Range("C13").Value = 0.2668690749977429
$ws.Range("D13").Value = 0.2098988909180264
$ws.Range("E13").Value = 0.1580149129163573
$ws.Range("F13").Value = 1.225763632141636
$ws.Range("H13").Value = 0.07973214163530429
$ws.Range("I13").Value = 0.5118691004029507
$ws.Range("J13").Value = 0.1600689384415617
$ws.Range("M13").Value = 0.4381801438148258
$ws.Range("N13").Value = 1.146131850495621
$ws.Range("O13").Value = 2.797742640936434

$ws.Range("B14").Value = 1.087309732006418
$ws.Range("C14").Value = 0.26244440410008
$ws.Range("D14").Value = 0.2094183144476602
$ws.Range("E14").Value = 0.157866982027695
$ws.Range("F14").Value = 1.224405216460767
$ws.Range("H14").Value = 0.07973214163530429
$ws.Range("I14").Value = 0.5125165105338532
$ws.Range("J14").Value = 0.1600800513221969
$ws.Range("M14").Value = 0.4333963976941746
$ws.Range("N14").Value = 1.147589288035718
$ws.Range("O14").Value = 2.795263039308026

$ws.Range("B15").Value = 1.076910186881832
$ws.Range("C15").Value = 0.2597336088314819
$ws.Range("D15").Value = 0.2091252448615393
$ws.Range("E15").Value = 0.1577776514281481
$ws.Range("F15").Value = 1.22358724210919
$ws.Range("H15").Value = 0.07973214163530429
$ws.Range("I15").Value = 0.5129197250847142
$ws.Range("J15").Value = 0.1600884058617069
$ws.Range("M15").Value = 0.4304679106346327
$ws.Range("N15").Value = 1.148490953553555
$ws.Range("O15").Value = 2.793779474379676

$ws.Range("B16").Value = 1.017309348628203
$ws.Range("C16").Value = 0.2441908209205508
$ws.Range("D16").Value = 0.2074657191061817
$ws.Range("E16").Value = 0.1572853910383287
$ws.Range("F16").Value = 1.219115758216986
$ws.Range("H16").Value = 0.07973214163530429
$ws.Range("I16").Value = 0.5153323392106444
$ws.Range("J16").Value = 0.1601600021235257
$ws.Range("M16").Value = 0.4137121839296896
$ws.Range("N16").Value = 1.153795134015027
$ws.Range("O16").Value = 2.785817956708058

$ws.Range("B17").Value = 0.9807425269557939
$ws.Range("C17").Value = 0.2346486369505669
$ws.Range("D17").Value = 0.2064655012970746
$ws.Range("E17").Value = 0.1570009649045261
$ws.Range("F17").Value = 1.216565348463746
$ws.Range("H17").Value = 0.07973214163530429
$ws.Range("I17").Value = 0.5169032870569943
$ws.Range("J17").Value = 0.1602250892587094
$ws.Range("M17").Value = 0.4034567935567637
$ws.Range("N17").Value = 1.15717144630581
$ws.Range("O17").Value = 2.781415534259054

$ws.Range("B18").Value = 0.9597093364213265
$ws.Range("C18").Value = 0.2291576765708498
$ws.Range("D18").Value = 0.2058967839038246
$ws.Range("E18").Value = 0.1568438312565377
$ws.Range("F18").Value = 1.21516922239816
$ws.Range("H18").Value = 0.07973214163530429
$ws.Range("I18").Value = 0.5178402202562253
$ws.Range("J18").Value = 0.1602703036015143
$ws.Range("M18").Value = 0.3975669880895083
$ws.Range("N18").Value = 1.159158395587156
$ws.Range("O18").Value = 2.779060279741913

$ws.Range("B19").Value = 0.9525877405669689
$ws.Range("C19").Value = 0.227298102813478
$ws.Range("D19").Value = 0.2057053584732671
$ws.Range("E19").Value = 0.1567917388644204
$ws.Range("F19").Value = 1.214708677037876
$ws.Range("H19").Value = 0.07973214163530429
$ws.Range("I19").Value = 0.5181631769639559
$ws.Range("J19").Value = 0.1602869482536207
$ws.Range("M19").Value = 0.3955743273085019
$ws.Range("N19").Value = 1.159838869798492
$ws.Range("O19").Value = 2.778293202352529

$ws.Range("B20").Value = 0.9846352323193628
$ws.Range("C20").Value = 0.2356646845793762
$ws.Range("D20").Value = 0.2065712955447481
$ws.Range("E20").Value = 0.1570305740615154
$ws.Range("F20").Value = 1.216829515842051
$ws.Range("H20").Value = 0.07973214163530429
$ws.Range("I20").Value = 0.5167326029726702
$ws.Range("J20").Value = 0.1602173556298894
$ws.Range("M20").Value = 0.40454758666732
$ws.Range("N20").Value = 1.156807376758387
$ws.Range("O20").Value = 2.78186586808576

$ws.Range("B21").Value = 1.092296488922329
$ws.Range("C21").Value = 0.2637441537028167
$ws.Range("D21").Value = 0.2095592011586263
$ws.Range("E21").Value = 0.1579101659275395
$ws.Range("F21").Value = 1.224801278718161
$ws.Range("H21").Value = 0.07973214163530429
$ws.Range("I21").Value = 0.5123249629452147
$ws.Range("J21").Value = 0.1600764646849626
$ws.Range("M21").Value = 0.4348011471578559
$ws.Range("N21").Value = 1.147159339667787
$ws.Range("O21").Value = 2.795984009016109

$ws.Range("B22").Value = 1.162642322198735
$ws.Range("C22").Value = 0.2820710478123374
$ws.Range("D22").Value = 0.2115700984425359
$ws.Range("E22").Value = 0.1585424003534968
$ws.Range("F22").Value = 1.230642082548655
$ws.Range("H22").Value = 0.07973214163530429
$ws.Range("I22").Value = 0.5097421791049541
$ws.Range("J22").Value = 0.1600536498056897
$ws.Range("M22").Value = 0.4546497892959493
$ws.Range("N22").Value = 1.14125408985025
$ws.Range("O22").Value = 2.806788968843136

$ws.Range("B23").Value = 1.125099410752682
$ws.Range("C23").Value = 0.2722919912889381
$ws.Range("D23").Value = 0.2104915491309214
$ws.Range("E23").Value = 0.1581997233863532
$ws.Range("F23").Value = 1.227467002841891
$ws.Range("H23").Value = 0.07973214163530429
$ws.Range("I23").Value = 0.5110933687503092
$ws.Range("J23").Value = 0.1600594885472759
$ws.Range("M23").Value = 0.4440493392536595
$ws.Range("N23").Value = 1.144369217799522
$ws.Range("O23").Value = 2.80087764940572

$ws.Range("B24").Value = 0.9828753726376931
$ws.Range("C24").Value = 0.2352053450989615
$ws.Range("D24").Value = 0.2065234462712198
$ws.Range("E24").Value = 0.1570171678661048
$ws.Range("F24").Value = 1.216709867296032
$ws.Range("H24").Value = 0.07973214163530429
$ws.Range("I24").Value = 0.5168096640494184
$ws.Range("J24").Value = 0.1602208277229877
$ws.Range("M24").Value = 0.4040544199037015
$ws.Range("N24").Value = 1.156971829736889
$ws.Range("O24").Value = 2.781661724716059

$ws.Range("B25").Value = 0.8294331560998671
$ws.Range("C25").Value = 0.1951049809720757
$ws.Range("D25").Value = 0.202495090591313
$ws.Range("E25").Value = 0.1559886840671147
$ws.Range("F25").Value = 1.207813176399526
$ws.Range("H25").Value = 0.07973214163530429
$ws.Range("I25").Value = 0.5242507989951015
$ws.Range("J25").Value = 0.1606922174407046
$ws.Range("M25").Value = 0.3612520267732222
$ws.Range("O25").Value = 2.767695776451376
